# Weekly fruit/vegetable price update: insert a new weekly record as row 18
# (pushing the existing rows 18-50 down to 19-51) and populate it with the
# latest price observation for Ají - Inferno - Primera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 18 - everything below shifts
# down by one (old row 18 becomes 19, ..., old row 50 becomes 51).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly observation.
$ws.Cells.Item(18, 1).Value  = 1
$ws.Cells.Item(18, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value  = 44498
$ws.Cells.Item(18, 5).Value  = 15
$ws.Cells.Item(18, 6).Value  = 100112021
$ws.Cells.Item(18, 7).Value  = "Ají"
$ws.Cells.Item(18, 8).Value  = "Inferno"
$ws.Cells.Item(18, 9).Value  = "Primera"
$ws.Cells.Item(18, 10).Value = 120
$ws.Cells.Item(18, 11).Value = 30000
$ws.Cells.Item(18, 12).Value = 31000
$ws.Cells.Item(18, 13).Value = 30500
$ws.Cells.Item(18, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(18, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(18, 16).Value = 2033
$ws.Cells.Item(18, 17).Value = 15
$ws.Cells.Item(18, 18).Value = "Hortaliza"
